$wb = $excel.ActiveWorkbook

# ---- Costs and Revenues ----
$wsCosts = $wb.Worksheets.Item("Costs and Revenues")
$wsCosts.Range("B2").Value = 76023.24100799997
$wsCosts.Range("D2").Value = 9272.28964544585
$wsCosts.Range("F2").Value = 16579.01275901571

# ---- Capacities ----
$wsCap = $wb.Worksheets.Item("Capacities")
$wsCap.Range("C2").Value = 1
$wsCap.Range("C3").Value = 103
$wsCap.Range("C4").Value = 161
$wsCap.Range("D4").Value = 0

# ---- PV Dispatch ----
$wsPV = $wb.Worksheets.Item("PV Dispatch")
$wsPV.Range("G2").Value = 20.6
$wsPV.Range("H2").Value = 41.2
$wsPV.Range("I2").Value = 51.5
$wsPV.Range("J2").Value = 61.8
$wsPV.Range("K2").Value = 72.09999999999999
$wsPV.Range("L2").Value = 82.40000000000001
$wsPV.Range("M2").Value = 92.7
$wsPV.Range("N2").Value = 103
$wsPV.Range("O2").Value = 92.7
$wsPV.Range("P2").Value = 82.40000000000001
$wsPV.Range("Q2").Value = 72.09999999999999
$wsPV.Range("R2").Value = 51.5
$wsPV.Range("S2").Value = 30.9
$wsPV.Range("T2").Value = 20.6

$wsPV.Range("K3").Value = 0
$wsPV.Range("L3").Value = 92.7
$wsPV.Range("M3").Value = 103
$wsPV.Range("N3").Value = 82.40000000000001
$wsPV.Range("O3").Value = 51.65717783899595
$wsPV.Range("P3").Value = 51.5
$wsPV.Range("Q3").Value = 51.5
$wsPV.Range("R3").Value = 30.9
$wsPV.Range("S3").Value = 20.6

$wsPV.Range("K4").Value = 41.2
$wsPV.Range("L4").Value = 72.09999999999999
$wsPV.Range("M4").Value = 23.4
$wsPV.Range("N4").Value = 0
$wsPV.Range("O4").Value = 59.98312417100291
$wsPV.Range("P4").Value = 41.2
$wsPV.Range("Q4").Value = 20.6
$wsPV.Range("R4").Value = 10.3

# ---- Battery Input ----
$wsBI = $wb.Worksheets.Item("Battery Input")
$wsBI.Range("G2").Value = 12.8
$wsBI.Range("H2").Value = 28.2
$wsBI.Range("I2").Value = 20.3
$wsBI.Range("J2").Value = 22.8
$wsBI.Range("K2").Value = 46.1
$wsBI.Range("L2").Value = 61.6
$wsBI.Range("M2").Value = 69.3
$wsBI.Range("N2").Value = 77
$wsBI.Range("O2").Value = 61.5
$wsBI.Range("P2").Value = 53.8
$wsBI.Range("Q2").Value = 46.1
$wsBI.Range("R2").Value = 17.7

$wsBI.Range("K3").Value = 0
$wsBI.Range("L3").Value = 92.7
$wsBI.Range("M3").Value = 79.59999999999999
$wsBI.Range("N3").Value = 56.4
$wsBI.Range("O3").Value = 51.65717783899595
$wsBI.Range("P3").Value = 22.9
$wsBI.Range("Q3").Value = 25.5
$wsBI.Range("R3").Value = 30.9

$wsBI.Range("K4").Value = 41.2
$wsBI.Range("L4").Value = 72.09999999999999
$wsBI.Range("M4").Value = 0
$wsBI.Range("N4").Value = 0
$wsBI.Range("O4").Value = 59.98312417100291
$wsBI.Range("P4").Value = 41.2
$wsBI.Range("Q4").Value = 20.6
$wsBI.Range("R4").Value = 10.3

# ---- Battery Output ----
$wsBO = $wb.Worksheets.Item("Battery Output")
$wsBO.Range("S2").Value = 1.007719999999879
$wsBO.Range("T2").Value = 31.4
$wsBO.Range("S3").Value = 21

# ---- State of Charge ----
$wsSOC = $wb.Worksheets.Item("State of Charge")
$wsSOC.Range("B2").Value = 187.8909090909091
$wsSOC.Range("C2").Value = 168.1939393939394
$wsSOC.Range("D2").Value = 155.0626262626263
$wsSOC.Range("E2").Value = 141.9313131313131
$wsSOC.Range("F2").Value = 128.8
$wsSOC.Range("G2").Value = 141.472
$wsSOC.Range("H2").Value = 169.39
$wsSOC.Range("I2").Value = 189.487
$wsSOC.Range("J2").Value = 212.059
$wsSOC.Range("K2").Value = 257.698
$wsSOC.Range("L2").Value = 318.682
$wsSOC.Range("M2").Value = 387.289
$wsSOC.Range("N2").Value = 463.519
$wsSOC.Range("O2").Value = 524.404
$wsSOC.Range("P2").Value = 577.6659999999999
$wsSOC.Range("Q2").Value = 623.3049999999999
$wsSOC.Range("R2").Value = 640.828
$wsSOC.Range("S2").Value = 639.8101010101011
$wsSOC.Range("T2").Value = 608.0929292929294
$wsSOC.Range("U2").Value = 489.9111111111111
$wsSOC.Range("V2").Value = 391.4262626262627
$wsSOC.Range("W2").Value = 312.6383838383838
$wsSOC.Range("X2").Value = 260.1131313131313
$wsSOC.Range("Y2").Value = 220.7191919191919

$wsSOC.Range("B3").Value = 181.3252525252525
$wsSOC.Range("C3").Value = 161.6282828282828
$wsSOC.Range("D3").Value = 148.4969696969697
$wsSOC.Range("E3").Value = 148.4969696969697
$wsSOC.Range("F3").Value = 148.4969696969697
$wsSOC.Range("G3").Value = 128.8
$wsSOC.Range("H3").Value = 128.8
$wsSOC.Range("I3").Value = 128.8
$wsSOC.Range("J3").Value = 128.8
$wsSOC.Range("K3").Value = 128.8
$wsSOC.Range("L3").Value = 220.573
$wsSOC.Range("M3").Value = 299.377
$wsSOC.Range("N3").Value = 355.213
$wsSOC.Range("O3").Value = 406.353606060606
$wsSOC.Range("P3").Value = 429.024606060606
$wsSOC.Range("Q3").Value = 454.269606060606
$wsSOC.Range("R3").Value = 484.860606060606
$wsSOC.Range("S3").Value = 463.6484848484848
$wsSOC.Range("T3").Value = 332.3353535353535
$wsSOC.Range("U3").Value = 332.3353535353535
$wsSOC.Range("V3").Value = 332.3353535353535
$wsSOC.Range("W3").Value = 253.5474747474748
$wsSOC.Range("X3").Value = 253.5474747474748
$wsSOC.Range("Y3").Value = 214.1535353535353

$wsSOC.Range("B4").Value = 168.1939393939394
$wsSOC.Range("C4").Value = 148.4969696969697
$wsSOC.Range("D4").Value = 148.4969696969697
$wsSOC.Range("E4").Value = 148.4969696969697
$wsSOC.Range("F4").Value = 148.4969696969697
$wsSOC.Range("G4").Value = 128.8
$wsSOC.Range("H4").Value = 128.8
$wsSOC.Range("I4").Value = 128.8
$wsSOC.Range("J4").Value = 128.8
$wsSOC.Range("K4").Value = 169.588
$wsSOC.Range("L4").Value = 240.967
$wsSOC.Range("M4").Value = 240.967
$wsSOC.Range("N4").Value = 240.967
$wsSOC.Range("O4").Value = 300.3502929292929
$wsSOC.Range("P4").Value = 341.1382929292929
$wsSOC.Range("Q4").Value = 361.5322929292929
$wsSOC.Range("R4").Value = 371.7292929292929
$wsSOC.Range("S4").Value = 371.7292929292929
$wsSOC.Range("T4").Value = 240.4161616161616
$wsSOC.Range("U4").Value = 240.4161616161616
$wsSOC.Range("V4").Value = 240.4161616161616
$wsSOC.Range("W4").Value = 240.4161616161616
$wsSOC.Range("X4").Value = 240.4161616161616
$wsSOC.Range("Y4").Value = 201.0222222222222

# ---- Feed in from Type 2 ----
$wsF2 = $wb.Worksheets.Item("Feed in from Type 2")
$wsF2.Range("S2").Value = 9.692280000000125
